$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.93

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.8

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.73

$ws.Range("B5").Value = 0.96
$ws.Range("C5").Value = 0.47

$ws.Range("B6").Value = 0.93
$ws.Range("C6").Value = 0.27

$ws.Range("B7").Value = 0.2
